$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Activate()

# --- Add new row 44: ammo_og-7b (rocket, EX type) ---
# Copy formats from the row above (row 43) so number formats / styles match,
# then fix up the two cells (A, C) that should keep the default (no) style.
$ws.Range("A43:K43").Copy()
$ws.Range("A44").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A44").ClearFormats()
$ws.Range("C44").ClearFormats()
$excel.CutCopyMode = 0

$ws.Range("A44").Value = "ammo_og-7b"
$ws.Range("B44").Value = "EX"
$ws.Range("C44").Value = 12490
$ws.Range("D44").Formula = "=C44/30"
$ws.Range("E44").Formula = "=K44/D44"
$ws.Range("F44").Formula = "=G44/D44*100"
$ws.Range("G44").Value = 0.37
$ws.Range("H44").Value = 3
$ws.Range("I44").Value = 3
$ws.Range("J44").Formula = "=I44*H44"
$ws.Range("K44").Formula = "=J44*Feuil2!`$B`$1"

# --- Restore the view state recorded in the workbook (scroll position / selection) ---
$null = $excel.Goto($ws.Range("A10"), $true)
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("D26").Select()
